$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.8499757026326253
$ws.Range("J4").Value = 0.4695572284139592
$ws.Range("K4").Value = 0.6273898113437535
$ws.Range("L4").Value = 3.108478537897936
